$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor-expressing-cell counts (1 -> 3) and all dependent NATMI
# specificity/weight metrics recomputed from the new single-cell counts,
# per "Natmi following Dr Hou advice".
$updates = @{
    "E2"=3; "G2"=19.95578266666667; "H2"=59.867348; "I2"=0.0117373419656925; "J2"=0.0117373419656925; "K2"=3; "M2"=153.5290173333333; "N2"=460.587052; "O2"=0.3172206968818489; "P2"=0.317220696881849; "Q2"=3063.791702930899; "R2"=27574.12532637809; "S2"=0.003723327797897544; "T2"=0.003723327797897545
    "E3"=3; "G3"=19.95578266666667; "H3"=59.867348; "I3"=0.0117373419656925; "J3"=0.0117373419656925; "K3"=3; "M3"=168.7997026666667; "N3"=506.3991080000001; "O3"=0.3487728915577651; "P3"=0.3487728915577651; "Q3"=3368.530180613955; "R3"=30316.77162552559; "S3"=0.004093666696576874; "T3"=0.004093666696576874
    "E4"=3; "G4"=19.95578266666667; "H4"=59.867348; "I4"=0.0117373419656925; "J4"=0.0117373419656925; "K4"=3; "M4"=68.09032333333333; "N4"=204.27097; "O4"=0.1406878008722904; "P4"=0.1406878008722904; "Q4"=1358.795694143062; "R4"=12229.16124728756; "S4"=0.001651300829239323; "T4"=0.001651300829239324
    "E5"=3; "G5"=19.95578266666667; "H5"=59.867348; "I5"=0.0117373419656925; "J5"=0.0117373419656925; "K5"=3; "M5"=93.562673; "N5"=280.688019; "O5"=0.1933186106880956; "P5"=0.1933186106880956; "Q5"=1867.116368100402; "R5"=16804.04731290361; "S5"=0.002269046641978755; "T5"=0.002269046641978755
    "E6"=3; "G6"=1637.343343333333; "H6"=4912.03003; "I6"=0.9630320723052701; "J6"=0.9630320723052702; "K6"=3; "M6"=153.5290173333333; "N6"=460.587052; "O6"=0.3172206968818489; "P6"=0.3172206968818489; "Q6"=251379.7145392412; "R6"=2262417.430853171; "S6"=0.3054937050962489; "T6"=0.305493705096249
    "E7"=3; "G7"=1637.343343333333; "H7"=4912.03003; "I7"=0.9630320723052701; "J7"=0.9630320723052702; "K7"=3; "M7"=168.7997026666667; "N7"=506.3991080000001; "O7"=0.3487728915577651; "P7"=0.3487728915577651; "Q7"=276383.0695179126; "R7"=2487447.625661213; "S7"=0.3358794805207758; "T7"=0.3358794805207758
    "E8"=3; "G8"=1637.343343333333; "H8"=4912.03003; "I8"=0.9630320723052701; "J8"=0.9630320723052702; "K8"=3; "M8"=68.09032333333333; "N8"=204.27097; "O8"=0.1406878008722904; "P8"=0.1406878008722904; "Q8"=111487.2376552477; "R8"=1003385.138897229; "S8"=0.135486864422113; "T8"=0.135486864422113
    "E9"=3; "G9"=1637.343343333333; "H9"=4912.03003; "I9"=0.9630320723052701; "J9"=0.9630320723052702; "K9"=3; "M9"=93.562673; "N9"=280.688019; "O9"=0.1933186106880956; "P9"=0.1933186106880956; "Q9"=153194.2198210234; "R9"=1378747.978389211; "S9"=0.1861720222661325; "T9"=0.1861720222661325
    "E10"=3; "G10"=17.50081933333334; "H10"=52.502458; "I10"=0.01029341242216722; "J10"=0.01029341242216722; "K10"=3; "M10"=153.5290173333333; "N10"=460.587052; "O10"=0.3172206968818489; "P10"=0.317220696881849; "Q10"=2686.883594774868; "R10"=24181.95235297382; "S10"=0.003265283461852165; "T10"=0.003265283461852166
    "E11"=3; "G11"=17.50081933333334; "H11"=52.502458; "I11"=0.01029341242216722; "J11"=0.01029341242216722; "K11"=3; "M11"=168.7997026666667; "N11"=506.3991080000001; "O11"=0.3487728915577651; "P11"=0.3487728915577651; "Q11"=2954.133099889719; "R11"=26587.19789900747; "S11"=0.003590063214475879; "T11"=0.00359006321447588
    "E12"=3; "G12"=17.50081933333334; "H12"=52.502458; "I12"=0.01029341242216722; "J12"=0.01029341242216722; "K12"=3; "M12"=68.09032333333333; "N12"=204.27097; "O12"=0.1406878008722904; "P12"=0.1406878008722904; "Q12"=1191.636447004918; "R12"=10724.72802304426; "S12"=0.001448157557146222; "T12"=0.001448157557146222
    "E13"=3; "G13"=17.50081933333334; "H13"=52.502458; "I13"=0.01029341242216722; "J13"=0.01029341242216722; "K13"=3; "M13"=93.562673; "N13"=280.688019; "O13"=0.1933186106880956; "P13"=0.1933186106880956; "Q13"=1637.423436516745; "R13"=14736.8109286507; "S13"=0.001989908188692952; "T13"=0.001989908188692952
    "E14"=3; "G14"=25.39612333333333; "H14"=76.18836999999999; "I14"=0.01493717330687017; "J14"=0.01493717330687017; "K14"=3; "M14"=153.5290173333333; "N14"=460.587052; "O14"=0.3172206968818489; "P14"=0.317220696881849; "Q14"=3899.041859442804; "R14"=35091.37673498524; "S14"=0.004738380525850306; "T14"=0.004738380525850307
    "E15"=3; "G15"=25.39612333333333; "H15"=76.18836999999999; "I15"=0.01493717330687017; "J15"=0.01493717330687017; "K15"=3; "M15"=168.7997026666667; "N15"=506.3991080000001; "O15"=0.3487728915577651; "P15"=0.3487728915577651; "Q15"=4286.858067552663; "R15"=38581.72260797396; "S15"=0.005209681125936573; "T15"=0.005209681125936573
    "E16"=3; "G16"=25.39612333333333; "H16"=76.18836999999999; "I16"=0.01493717330687017; "J16"=0.01493717330687017; "K16"=3; "M16"=68.09032333333333; "N16"=204.27097; "O16"=0.1406878008722904; "P16"=0.1406878008722904; "Q16"=1729.230249179878; "R16"=15563.0722426189; "S16"=0.002101478063791842; "T16"=0.002101478063791842
    "E17"=3; "G17"=25.39612333333333; "H17"=76.18836999999999; "I17"=0.01493717330687017; "J17"=0.01493717330687017; "K17"=3; "M17"=93.562673; "N17"=280.688019; "O17"=0.1933186106880956; "P17"=0.1933186106880956; "Q17"=2376.129182904337; "R17"=21385.16264613903; "S17"=0.002887633591291448; "T17"=0.002887633591291448
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
